# Auto-generated edit script: refresh crypto price/volume data (Mon Jan 30 21:30:26 UTC 2023 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.30"
$ws.Range("E2").Value = "'-3.63%"
$ws.Range("D3").Value = "'37.49"
$ws.Range("E4").Value = "'-0.97%"
$ws.Range("D5").Value = "'0.07725"
$ws.Range("E5").Value = "'-5.77%"
$ws.Range("D6").Value = "'4.351"
$ws.Range("E6").Value = "'0.28%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.890"
$ws.Range("E7").Value = "'-7.18%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.207"
$ws.Range("E8").Value = "'-1.74%"
$ws.Range("E9").Value = "'-1.13%"
$ws.Range("D10").Value = "'0.9168"
$ws.Range("E10").Value = "'-2.58%"
$ws.Range("D11").Value = "'0.1189"
$ws.Range("E11").Value = "'-12.21%"
$ws.Range("D12").Value = "'0.1885"
$ws.Range("E12").Value = "'-5.05%"
$ws.Range("D13").Value = "'0.08702"
$ws.Range("E13").Value = "'-4.58%"
$ws.Range("D14").Value = "'0.03406"
$ws.Range("E14").Value = "'-2.57%"
$ws.Range("D15").Value = "'0.09704"
$ws.Range("E15").Value = "'-0.82%"
$ws.Range("D16").Value = "'0.001368"
$ws.Range("E16").Value = "'-2.81%"
$ws.Range("D17").Value = "'0.005946"
$ws.Range("E17").Value = "'-1.89%"
$ws.Range("D18").Value = "'3.584"
$ws.Range("E18").Value = "'-2.87%"
$ws.Range("D19").Value = "'0.3408"
$ws.Range("D20").Value = "'0.1275"
$ws.Range("E20").Value = "'-2.91%"
$ws.Range("D21").Value = "'5.019"
$ws.Range("E21").Value = "'1.25%"
$ws.Range("D22").Value = "'0.2596"
$ws.Range("E22").Value = "'5.87%"
$ws.Range("E23").Value = "'5,161.13%"
$ws.Range("E24").Value = "'-0.82%"
$ws.Range("E25").Value = "'-1.37%"
$ws.Range("D26").Value = "'0.004544"
$ws.Range("E26").Value = "'-5.23%"
$ws.Range("E27").Value = "'3.88%"
$ws.Range("D39").Value = "'0.02213"
$ws.Range("E39").Value = "'-2.42%"
$ws.Range("D40").Value = "'0.04917"
$ws.Range("E40").Value = "'-5.36%"
$ws.Range("D41").Value = "'0.007552"
$ws.Range("E41").Value = "'-2.61%"
$ws.Range("D42").Value = "'0.009903"
$ws.Range("E42").Value = "'-0.02%"
$ws.Range("D43").Value = "'0.1336"
$ws.Range("E43").Value = "'-4.67%"
$ws.Range("D44").Value = "'0.002064"
$ws.Range("E44").Value = "'0.87%"
$ws.Range("D45").Value = "'0.008806"
$ws.Range("E45").Value = "'-3.32%"
$ws.Range("D46").Value = "'0.00006558"
$ws.Range("E46").Value = "'-0.65%"
$ws.Range("E47").Value = "'0.05%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.002999"
$ws.Range("E48").Value = "'1.70%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.001302"
$ws.Range("E49").Value = "'-23.03%"
$ws.Range("E50").Value = "'0.05%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.05%"

# The apostrophe prefix above leaves a quotePrefix style on each touched cell;
# clear formats on the touched numeric range so styles match the original (unstyled) cells
# while keeping the values as literal text.
$ws.Range("D2:E51").ClearFormats()
